$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "'111011101100101101111011011001"
$ws.Range("D2").Value = 0.8701016422874415
$ws.Range("E2").Value = 0.0000990339324305849
$ws.Range("F2").Value = 0.3868239948902124

$ws.Range("C3").Value = "'111011101100101101111011011101"
$ws.Range("D3").Value = 0.8701016492372863
$ws.Range("E3").Value = 0.07563701921514805
$ws.Range("F3").Value = 0.628301577707689

$ws.Range("C4").Value = "'111011101100101101111011011101"
$ws.Range("D4").Value = 0.8701016492372863
$ws.Range("E4").Value = 0.1825878186348867
$ws.Range("F4").Value = 0.7568158355382901

$ws.Range("C5").Value = "'111011101100101101111011011101"
$ws.Range("D5").Value = 0.8701016492372863
$ws.Range("E5").Value = 0.6525310182886788
$ws.Range("F5").Value = 0.8240172570352573

$ws.Range("C6").Value = "'111011101100101101111011011101"
$ws.Range("D6").Value = 0.8701016492372863
$ws.Range("E6").Value = 0.8100297129457867
$ws.Range("F6").Value = 0.8512231405743359

$ws.Range("C7").Value = "'111011101100101111010101000101"
$ws.Range("D7").Value = 0.8701116153434306
$ws.Range("E7").Value = 0.810039342281386
$ws.Range("F7").Value = 0.8580806827317223

$ws.Range("C8").Value = "'111011101100101111010101000101"
$ws.Range("D8").Value = 0.8701116153434306
$ws.Range("E8").Value = 0.8672573210669913
$ws.Range("F8").Value = 0.8698192096414857

$ws.Range("C9").Value = "'111011101101101111010101000101"
$ws.Range("D9").Value = 0.870567142592784
$ws.Range("E9").Value = 0.8701016492372863
$ws.Range("F9").Value = 0.8701521850152938

$ws.Range("C10").Value = "'111011101101101111010101000101"
$ws.Range("D10").Value = 0.870567142592784
$ws.Range("E10").Value = 0.8701007596573647
$ws.Range("F10").Value = 0.8702900356003564

$ws.Range("C11").Value = "'111011101101101111010101001101"
$ws.Range("D11").Value = 0.8705671564961914
$ws.Range("E11").Value = 0.870100745757682
$ws.Range("F11").Value = 0.8703372926193518

$ws.Range("C12").Value = "'111011101101101111010101001101"
$ws.Range("D12").Value = 0.8705671564961914
$ws.Range("E12").Value = 0.4662127472782915
$ws.Range("F12").Value = 0.8299929599281141

$ws.Range("C13").Value = "'111011101101101111010101001101"
$ws.Range("D13").Value = 0.8705671564961914
$ws.Range("E13").Value = 0.4662127472782915
$ws.Range("F13").Value = 0.7896497155857919

$ws.Range("C14").Value = "'111111101101101111010101000101"
$ws.Range("D14").Value = 0.9911036272236304
$ws.Range("E14").Value = 0.8701016075382174
$ws.Range("F14").Value = 0.8825276868256369

$ws.Range("C15").Value = "'111111101101101111010101000101"
$ws.Range("D15").Value = 0.9911036272236304
$ws.Range("E15").Value = 0.8701016075382174
$ws.Range("F15").Value = 0.8825248394310824

$ws.Range("C16").Value = "'111111101101101111010101000101"
$ws.Range("D16").Value = 0.9911036272236304
$ws.Range("E16").Value = 0.8701016075382174
$ws.Range("F16").Value = 0.882523436583097

$ws.Range("C17").Value = "'111111101101101111010101000101"
$ws.Range("D17").Value = 0.9911036272236304
$ws.Range("E17").Value = 0.870567142592784
$ws.Range("F17").Value = 0.8946744617644065

$ws.Range("C18").Value = "'111111101101101111010101010101"
$ws.Range("D18").Value = 0.9911036568930905
$ws.Range("E18").Value = 0.870567142592784
$ws.Range("F18").Value = 0.9067280909489839

$ws.Range("C19").Value = "'111111101101101111010101010101"
$ws.Range("D19").Value = 0.9911036568930905
$ws.Range("E19").Value = 0.8705671356410803
$ws.Range("F19").Value = 0.9308353871799827

$ws.Range("C20").Value = "'111111101101101111010101010101"
$ws.Range("D20").Value = 0.9911036568930905
$ws.Range("E20").Value = 0.870567142592784
$ws.Range("F20").Value = 0.9549426959273701

$ws.Range("C21").Value = "'111111101101101111010101010111"
$ws.Range("D21").Value = 0.9911036606017731
$ws.Range("E21").Value = 0.870567142592784
$ws.Range("F21").Value = 0.9549427041779429
